$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 22 and 23 had their match data (columns F, H..V) swapped upstream.
#    Columns A-D, E and G stay put per row position.
#    Swap via Value2 (preserves numeric vs. string typing) cell-by-cell.
# ---------------------------------------------------------------------------
$swapCols = @(6,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22)  # F, H..V (skip G)
foreach ($c in $swapCols) {
    $v1 = $ws.Cells.Item(22, $c).Value2
    $v2 = $ws.Cells.Item(23, $c).Value2
    $ws.Cells.Item(22, $c).Value = $v2
    $ws.Cells.Item(23, $c).Value = $v1
}

# ---------------------------------------------------------------------------
# 2) Append three new match rows (99, 100, 101) at the bottom of the sheet.
# ---------------------------------------------------------------------------
$newRows = @(
    @{
        Row = 99; Indice = 98; E = 45234.66666666666
        F = "RWDM"; G = 1; H = "Kortrijk"; I = 1
        J = 1.78; K = "29/10/2023 19:43"; L = 2.05; M = "04/11/2023 15:52"
        N = 4.03; O = "29/10/2023 19:43"; P = 3.75; Q = "04/11/2023 15:56"
        R = 3.92; S = "29/10/2023 19:43"; T = 3.58; U = "04/11/2023 15:42"
        V = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/rwd-molenbeek-kortrijk/bH3MdCA1/"
    }
    @{
        Row = 100; Indice = 99; E = 45234.76041666666
        F = "Leuven"; G = 0; H = "Westerlo"; I = 2
        J = 2.19; K = "29/10/2023 16:12"; L = 2; M = "04/11/2023 18:13"
        N = 3.69; O = "29/10/2023 16:12"; P = 3.97; Q = "04/11/2023 18:13"
        R = 3.01; S = "29/10/2023 16:12"; T = 3.54; U = "04/11/2023 18:13"
        V = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/leuven-westerlo/2TrxIMtA/"
    }
    @{
        Row = 101; Indice = 100; E = 45234.86458333334
        F = "Antwerp"; G = 3; H = "Genk"; I = 2
        J = 2.08; K = "29/10/2023 19:43"; L = 2.22; M = "04/11/2023 20:16"
        N = 3.74; O = "29/10/2023 19:43"; P = 3.75; Q = "04/11/2023 20:16"
        R = 3.18; S = "29/10/2023 19:43"; T = 3.16; U = "04/11/2023 20:16"
        V = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/antwerp-genk/8SA9aEQr/"
    }
)

foreach ($row in $newRows) {
    $r = $row.Row

    # Column A - reuse the bold/bordered style used throughout the sheet.
    $ws.Cells.Item(98, 1).Copy($ws.Cells.Item($r, 1))
    $ws.Cells.Item($r, 1).Value = $row.Indice

    $ws.Cells.Item($r, 2).Value = "belgium"
    $ws.Cells.Item($r, 3).Value = "jupiler-pro-league"
    $ws.Cells.Item($r, 4).Value = "2023-2024"

    # Column E - reuse the custom datetime number format.
    $ws.Cells.Item(98, 5).Copy($ws.Cells.Item($r, 5))
    $ws.Cells.Item($r, 5).Value = $row.E

    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $ws.Cells.Item($r, 21).Value = $row.U
    $ws.Cells.Item($r, 22).Value = $row.V
}
